# Initial Sync of Parallel Development Code Base
#
# Adds two new third-level ("Add new Monitoring Unit Id" / "Add new Entry")
# bullet items under the existing "Field Fortifications" bullet (ilvl=2,
# numId=4 / the "Data" > "Studies" > "Monitoring Units" list), right before
# the trailing blank paragraph at the end of the document.

$d = $word.ActiveDocument

# Locate the paragraph that anchors the insertion point: the last bullet
# ("Field Fortifications") under the "Data" heading's list.
$anchor = $d.Content
$found = $anchor.Find.Execute("Field Fortifications", $true, $false, $false,
                               $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the 'Field Fortifications' paragraph to anchor the new list items on."
}

# Collapse to right after the found text, i.e. just before that paragraph's
# own end-of-paragraph mark, so the new paragraphs get inserted immediately
# after it (and before the document's trailing empty paragraph).
$insPoint = $d.Range($anchor.End, $anchor.End)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newPara1 = "<w:p $wNs>" +
              "<w:pPr>" +
                "<w:pStyle w:val=`"ListParagraph`"/>" +
                "<w:numPr><w:ilvl w:val=`"2`"/><w:numId w:val=`"4`"/></w:numPr>" +
              "</w:pPr>" +
              "<w:r><w:t>Add new Monitoring Unit Id</w:t></w:r>" +
            "</w:p>"

$newPara2 = "<w:p $wNs>" +
              "<w:pPr>" +
                "<w:pStyle w:val=`"ListParagraph`"/>" +
                "<w:numPr><w:ilvl w:val=`"2`"/><w:numId w:val=`"4`"/></w:numPr>" +
              "</w:pPr>" +
              "<w:r><w:lastRenderedPageBreak/><w:t>Add new Entry</w:t></w:r>" +
            "</w:p>"

$insPoint.InsertXML($newPara1 + $newPara2)

Write-Output "Inserted 'Add new Monitoring Unit Id' and 'Add new Entry' bullets after 'Field Fortifications'."
